$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the out-of-stock capacitor assortment with the new Wandefol product
$newUrl = "https://www.amazon.fr/Wandefol-Condensateur-Electrolytique-Rangement-Valeurs/dp/B07Q6PNB6H/"

$ws.Range("D7").Hyperlinks.Delete()
$ws.Range("D7").Value = $newUrl
$ws.Range("E7").Value = 15.95
$ws.Hyperlinks.Add($ws.Range("D7"), $newUrl)

# Add a Total row summing the Price column
$ws.Range("E20").Value = "Total"
$ws.Range("F20").Formula = "=SUM(F3:F18)"

$ws.Range("G20").Select()
